# Fix version of PSYIVX (UK EPU Index) data:
#  - update May 2021 value
#  - append Jun-Dec 2021 rows (previously missing)
#  - the trailing "Source: ..." note row shifts down accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the existing May 2021 (row 282) value.
$ws.Range("C282").Value = 98.111930847167969

# 2) Insert 7 new blank rows right after row 282 (before the "Source:" row)
#    so the trailing note row is pushed down from 283 to 290.
$ws.Range("A283:A289").EntireRow.Insert()

# 3) Seed each new row by copying row 282 (year/month/value columns),
#    which keeps column A as the shared "2021" text value, then
#    overwrite the month (B) and index value (C) for each month.
$newRows = @(
    @{ Row = 283; Month = 6;  Value = 90.334930419921875 },
    @{ Row = 284; Month = 7;  Value = 88.283767700195313 },
    @{ Row = 285; Month = 8;  Value = 79.109107971191406 },
    @{ Row = 286; Month = 9;  Value = 78.700927734375 },
    @{ Row = 287; Month = 10; Value = 103.98979187011719 },
    @{ Row = 288; Month = 11; Value = 75.317497253417969 },
    @{ Row = 289; Month = 12; Value = 164.3631591796875 }
)

foreach ($r in $newRows) {
    $ws.Range("A282:C282").Copy($ws.Range("A" + $r.Row + ":C" + $r.Row))
    $ws.Range("B" + $r.Row).Value = $r.Month
    $ws.Range("C" + $r.Row).Value = $r.Value
}
